# Update "想去人数" (want-to-go count) values in column F for the
# "展览" and "全部类型" worksheets, matching the regenerated data output.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 159
    3  = 421
    4  = 12281
    5  = 1269
    11 = 446
    17 = 2962
    18 = 89
    20 = 11
    22 = 21
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
